$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.993.81"
$ws.Range("E2").Value = "  -0.52%  "
$ws.Range("D3").Value = "3.831.57"
$ws.Range("E3").Value = "  +0.78%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = $ws.Range("B4").Style
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "701.68"
$ws.Range("D5").Style = $ws.Range("B5").Style
$ws.Range("E5").Value = "  +0.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.78"
$ws.Range("D6").Style = $ws.Range("B6").Style
$ws.Range("E6").Value = "  -1.26%  "
$ws.Range("D7").Value = "3.830.36"
$ws.Range("E7").Value = "  +0.81%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.526"
$ws.Range("D9").Style = $ws.Range("B9").Style
$ws.Range("E9").Value = "  -0.51%  "
$ws.Range("E10").Value = "  -0.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.34"
$ws.Range("D11").Style = $ws.Range("B11").Style
$ws.Range("E11").Value = "  -0.82%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.459"
$ws.Range("D12").Style = $ws.Range("B12").Style
$ws.Range("E12").Value = "  -0.42%  "
$ws.Range("E13").Value = "  -1.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.70"
$ws.Range("D14").Style = $ws.Range("B14").Style
$ws.Range("E14").Value = "  +0.35%  "
$ws.Range("D15").Value = "4.472.00"
$ws.Range("E15").Value = "  +0.57%  "
$ws.Range("D16").Value = "3.778.66"
$ws.Range("E16").Value = "  -0.72%  "
$ws.Range("D17").Value = "71.016.54"
$ws.Range("E17").Value = "  -0.47%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.23"
$ws.Range("D18").Style = $ws.Range("B18").Style
$ws.Range("E18").Value = "  -0.12%  "
$ws.Range("E19").Value = "  +0.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.38"
$ws.Range("D20").Style = $ws.Range("B20").Style
$ws.Range("E20").Value = "  -1.94%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "496.44"
$ws.Range("D21").Style = $ws.Range("B21").Style
$ws.Range("E21").Value = "  +2.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.69"
$ws.Range("D22").Style = $ws.Range("B22").Style
$ws.Range("E22").Value = "  -3.93%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.733"
$ws.Range("D23").Style = $ws.Range("B23").Style
$ws.Range("E23").Value = "  +2.40%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.15"
$ws.Range("D24").Style = $ws.Range("B24").Style
$ws.Range("E24").Value = "  +0.77%  "
$ws.Range("E25").Value = "  +1.22%  "
$ws.Range("E26").Value = "  +0.80%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.15"
$ws.Range("D27").Style = $ws.Range("B27").Style
$ws.Range("E27").Value = "  -1.40%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.09"
$ws.Range("D28").Style = $ws.Range("B28").Style
$ws.Range("E28").Value = "  -3.24%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.07"
$ws.Range("D30").Style = $ws.Range("B30").Style
$ws.Range("E30").Value = "  -1.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.45"
$ws.Range("D31").Style = $ws.Range("B31").Style
$ws.Range("E31").Value = "  -1.69%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.23"
$ws.Range("D32").Style = $ws.Range("B32").Style
$ws.Range("E32").Value = "  -3.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.39"
$ws.Range("D33").Style = $ws.Range("B33").Style
$ws.Range("E33").Value = "  -0.65%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.178"
$ws.Range("D34").Style = $ws.Range("B34").Style
$ws.Range("E34").Value = "  -3.28%  "
$ws.Range("B35").Value = "Aptos"
$ws.Range("C35").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.20"
$ws.Range("D35").Style = $ws.Range("B35").Style
$ws.Range("E35").Value = "  -0.87%  "
$ws.Range("B36").Value = "RenzoRestakedETH"
$ws.Range("C36").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D36").Value = "3.793.24"
$ws.Range("E36").Value = "  +1.03%  "
$ws.Range("B37").Value = "Binance-PegBSC-USD"
$ws.Range("C37").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").Style = $ws.Range("B37").Style
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("E38").Value = "  -0.81%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.32"
$ws.Range("D39").Style = $ws.Range("B39").Style
$ws.Range("E39").Value = "  -2.01%  "
$ws.Range("E40").Value = "  +3.97%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.97"
$ws.Range("D41").Style = $ws.Range("B41").Style
$ws.Range("E41").Value = "  -0.95%  "
$ws.Range("E42").Value = "  -2.72%  "
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("E44").Value = "  -0.15%  "
$ws.Range("E45").Value = "  +3.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "163.34"
$ws.Range("D46").Style = $ws.Range("B46").Style
$ws.Range("E46").Value = "  -0.83%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "432.87"
$ws.Range("D47").Style = $ws.Range("B47").Style
$ws.Range("E47").Value = "  +3.91%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "49.00"
$ws.Range("D48").Style = $ws.Range("B48").Style
$ws.Range("E48").Value = "  +0.82%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.72"
$ws.Range("D49").Style = $ws.Range("B49").Style
$ws.Range("E49").Value = "  +0.65%  "
$ws.Range("E50").Value = "  -0.08%  "
$ws.Range("E51").Value = "  -2.16%  "
